$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) "PSY" -> "PSYCHIC" : the Pokemon type label used by Mewtwo (D35) and Mew (D36)
$ws.Range("D35").Value = "PSYCHIC"
$ws.Range("D36").Value = "PSYCHIC"

# 2) The custom "boolean" number format (applied to the TRUE()/FALSE() helper
#    columns C and F) changes its display text from the French "VRAI/FAUX"
#    style to the "BOOL"E"AN" style already used elsewhere in the workbook.
$ws.Range("C1:C36").NumberFormat = """BOOL""E""AN"""
$ws.Range("F1:F36").NumberFormat = """BOOL""E""AN"""

# 3) Scroll the sheet view so row 13 is the top-left visible row, and move
#    the active selection to D36.
$ws.Range("D36").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
